$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.752.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.38%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.471.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.10%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.99"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.43%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.68"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.84%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.59%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.10%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.541"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.40%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.13"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.07%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0815"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.92%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.71%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.23"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.71%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.12"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.73%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.855.08"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.90%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.514.02"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.38%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.36%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.612.03"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.50%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.60"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.71%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.44%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.96%  "

# Row 23
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.74"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.03%  "

# Row 24
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.71%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.69%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.01%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.25"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.40%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.81"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.43%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.05"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.11%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.89%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.61"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.93%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.60"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.55%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.75%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.07%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.47%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.55%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.30%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.34%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "123.22"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.21%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.74%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.84%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.58"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.48%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0293"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.85%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.975.80"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.67%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.15%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.86%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.40%  "

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.96"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.41%  "

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.35"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +16.25%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.41"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.93%  "
